# Corrected excel sheets for application fix issues
#
# - Updates a couple of figures on the "Summary" sheet (row 3 and row 4).
# - Touches G2 on "Summary" so the sheet's used range grows to column G.
# - Makes "Summary" the active sheet/tab with D3 selected (it previously had
#   A7:XFD14 selected while "Transactions" was the tab shown on open).

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("Summary")

# Corrected figures
$summary.Range("A3").Value = 297.55
$summary.Range("E3").Value = 297.55

$summary.Range("A4").Value = 150
$summary.Range("E4").Value = 100

# Extend the sheet's used range out to column G (adds a blank G2 cell).
$summary.Range("G2").Style = "Normal"

# Make "Summary" the selected/active sheet, with D3 as the active cell.
$summary.Activate()
$summary.Range("D3").Select()
